$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 currently holds "Chlorophyll A (92nd Percentile)" for 2017-2021.
# This row is removed entirely, which shifts the old rows 8 (MCI) and
# 9 (QMCI) up to become rows 7 and 8, and the sheet shrinks from 9 to 8
# data rows (dimension A1:U9 -> A1:U8).
$ws.Rows.Item(7).Delete()
